# Apply updated "dSF" (final delta) values in column F of Sheet1.
# These reflect the repulled/pushed data and recalculated mean, per commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F6"  = -2
    "F7"  = 1
    "F8"  = -4
    "F13" = 6
    "F17" = -5
    "F18" = -3
    "F20" = -7
    "F21" = -12
    "F22" = -2
    "F23" = -11
    "F28" = -7
    "F29" = -4
    "F30" = 2
    "F36" = 1
    "F38" = -9
    "F40" = -1
    "F44" = -1
    "F49" = 3
    "F50" = -5
    "F53" = -10
    "F55" = -3
    "F59" = -4
    "F62" = -3
    "F63" = -2
    "F64" = 2
    "F66" = 7
    "F69" = 0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
